# Update forecast values (Working Multiple ASINs)
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: columns B (Prophet Forecast) and D (yhat_upper) ---

# Row 2
$wsForecast.Range("D2").Value = 234.0890616412142

# Row 3
$wsForecast.Range("B3").Value = 203
$wsForecast.Range("D3").Value = 255.5624876177054

# Row 4
$wsForecast.Range("D4").Value = 228.1469955562823

# Row 5
$wsForecast.Range("B5").Value = 125
$wsForecast.Range("D5").Value = 171.0383102333376

# Row 6
$wsForecast.Range("B6").Value = 84
$wsForecast.Range("D6").Value = 129.1274725553662

# Row 7
$wsForecast.Range("B7").Value = 71
$wsForecast.Range("D7").Value = 121.9366143077742

# Row 8
$wsForecast.Range("D8").Value = 122.4787996547277

# Row 9
$wsForecast.Range("B9").Value = 86
$wsForecast.Range("D9").Value = 138.1457550133834

# Row 10
$wsForecast.Range("B10").Value = 91
$wsForecast.Range("D10").Value = 142.975817229278

# Row 11
$wsForecast.Range("B11").Value = 85
$wsForecast.Range("D11").Value = 130.5190494994928

# Row 12
$wsForecast.Range("B12").Value = 78
$wsForecast.Range("D12").Value = 129.8849263716982

# Row 13
$wsForecast.Range("B13").Value = 64
$wsForecast.Range("D13").Value = 116.591158556338

# Row 14
$wsForecast.Range("D14").Value = 100.7306759217145

# Row 15
$wsForecast.Range("B15").Value = 48
$wsForecast.Range("D15").Value = 98.98452153795701

# Row 16
$wsForecast.Range("B16").Value = 53
$wsForecast.Range("D16").Value = 101.3594400008375

# Row 17
$wsForecast.Range("D17").Value = 110.3852022442807

# Row 18
$wsForecast.Range("B18").Value = 64
$wsForecast.Range("D18").Value = 114.247613863925

# Row 19
$wsForecast.Range("B19").Value = 58
$wsForecast.Range("D19").Value = 107.9919638190198

# Row 20
$wsForecast.Range("D20").Value = 104.796324783441

# Row 21
$wsForecast.Range("B21").Value = 69
$wsForecast.Range("D21").Value = 119.8240781494048

# --- Summary sheet: updated derived metrics ---
# Leading apostrophe forces these numeric-looking values to be stored as
# text (matching the source file's inlineStr cells); resetting the style
# back to Normal afterwards drops the quote-prefix formatting that Excel
# would otherwise apply, so no stray cell style is introduced.
$wsSummary.Range("B10").Value = "'1006"
$wsSummary.Range("B10").Style = "Normal"

$wsSummary.Range("B11").Value = "'691"
$wsSummary.Range("B11").Style = "Normal"

$wsSummary.Range("B12").Value = "'203"
$wsSummary.Range("B12").Style = "Normal"

$wsSummary.Range("B14").Value = "'48"
$wsSummary.Range("B14").Style = "Normal"
